$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.278.02"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.895.53"

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.94%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.343"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "50.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0709"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "2.171.82"
$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.693"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.72%  "

$ws.Range("D16").Value = "1.876.99"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("D18").Value = "35.316.91"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("E20").Value = "  +2.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +33.63%  "

$ws.Range("E26").Value = "  +1.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "

$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.928"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.94%  "

$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("E38").Value = "  +2.74%  "

$ws.Range("E39").Value = "  +3.60%  "

$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("E41").Value = "  +14.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "89.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").Value = "1.337.40"

$ws.Range("E45").Value = "  +2.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +38.24%  "

$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -14.74%  "
